# Adicionando tabela com escolhas para o jogador com os maiores scores
# para o esquema tatico 5-3-2
#
# Row 16 used to hold the "MEDIA 4-4-2 - Time de Maior Score (media)" line;
# it is repurposed here to hold the new "SCORE 5-3-2 - Time de Maior Score
# (media)" line, with refreshed per-round data and a percentage formula
# that now compares against the 5-3-2 scheme's total (row 8) instead of
# the 4-4-2 scheme's total (row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Re-label row 16.
$ws.Range("B16").Value = "SCORE 5-3-2 - Time de Maior Score (media)"

# 2) Refresh the 38 round scores (columns C..AN) for row 16.
$values = @(
    132.77999999999997, 124.44, 109.13000000000002, 128.43,
    113.92999999999999, 130.73000000000002, 97.039999999999992,
    115.78000000000002, 85.1, 125.75000000000001, 84.160000000000011,
    85.51, 77.839999999999989, 81.289999999999992, 67.7,
    89.12, 85.990000000000009, 93.09, 80.930000000000007,
    40.590000000000003, 37.36, 59.3, 50.650000000000006,
    75.489999999999995, 69.490000000000009, 51.45, 37.119999999999997,
    48.74, 28.54, 45.599999999999994, 72.319999999999993,
    42.65, 25.02, 52.960000000000008, 55.16,
    43.45, 108.56000000000003, 72.170000000000016
)
for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item(16, 3 + $i).Value = $values[$i]
}

# 3) The % formula on row 16 now compares against the 5-3-2 total (AO8)
#    instead of the 4-4-2 total (AO3).
$ws.Range("AQ16").Formula = '=(AO16*100)/$AO$8'

# 4) Both percentage cells (row 15 and row 16) get a "0.00" number format
#    with centered horizontal/vertical alignment.
$ws.Range("AQ15:AQ16").NumberFormat = "0.00"
$ws.Range("AQ15:AQ16").HorizontalAlignment = -4108
$ws.Range("AQ15:AQ16").VerticalAlignment = -4108

# 5) Leave the cursor where the author ended up after the edit.
$ws.Range("B17").Select() | Out-Null
